$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1490
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1490
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1490
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -1840
$ws.Range("H33").Value = 651.2778
$ws.Range("I33").Value = 384.3
$ws.Range("J33").Value = 985
$ws.Range("K33").Value = 384.3
$ws.Range("L33").Value = 985
$ws.Range("M33").Value = -155.3
$ws.Range("N33").Value = -1443
$ws.Range("H64").Value = 5115.5
$ws.Range("J64").Value = 5115.5
$ws.Range("L64").Value = 5115.5
$ws.Range("N64").Value = -5611.5
$ws.Range("H67").Value = 5115.5
$ws.Range("J67").Value = 5115.5
$ws.Range("L67").Value = 5115.5
$ws.Range("N67").Value = -6831.5
$ws.Range("H113").Value = 5140.4287
$ws.Range("I113").Value = 5498.8335
$ws.Range("K113").Value = 5498.8335
$ws.Range("M113").Value = -2244.8335
$ws.Range("H132").Value = 61750
$ws.Range("I132").Value = 2618
$ws.Range("J132").Value = 337699.34
$ws.Range("K132").Value = 7854
$ws.Range("L132").Value = 1013098.02
$ws.Range("M132").Value = -5324
$ws.Range("N132").Value = -1018158.02
$ws.Range("H137").Value = 3765.95
$ws.Range("I137").Value = 3140.1667
$ws.Range("K137").Value = 9420.500100000001
$ws.Range("M137").Value = -6870.500100000001
$ws.Range("H138").Value = 2228.1853
$ws.Range("J138").Value = 2168.9167
$ws.Range("L138").Value = 6506.750100000001
$ws.Range("N138").Value = -16786.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5290.7925
$ws.Range("I32").Value = 1798.5555
$ws.Range("K32").Value = 1798.5555
$ws.Range("M32").Value = -1511.5555
$ws.Range("H61").Value = 4468.1055
$ws.Range("I61").Value = 4159.1816
$ws.Range("K61").Value = 4159.1816
$ws.Range("M61").Value = -3947.1816
$ws.Range("H62").Value = 49999.5
$ws.Range("J62").Value = 49999.5
$ws.Range("L62").Value = 49999.5
$ws.Range("N62").Value = -51247.5
$ws.Range("H63").Value = 8889.666999999999
$ws.Range("I63").Value = 7207.143
$ws.Range("J63").Value = 10361.875
$ws.Range("K63").Value = 7207.143
$ws.Range("L63").Value = 10361.875
$ws.Range("M63").Value = -6521.143
$ws.Range("N63").Value = -11733.875
$ws.Range("H65").Value = 49999.5
$ws.Range("J65").Value = 49999.5
$ws.Range("L65").Value = 149998.5
$ws.Range("N65").Value = -156238.5
$ws.Range("H66").Value = 8889.666999999999
$ws.Range("I66").Value = 7207.143
$ws.Range("J66").Value = 10361.875
$ws.Range("K66").Value = 36035.715
$ws.Range("L66").Value = 51809.375
$ws.Range("M66").Value = -32603.715
$ws.Range("N66").Value = -58673.375
$ws.Range("H132").Value = 2070.7
$ws.Range("I132").Value = 1320.4736
$ws.Range("K132").Value = 3961.4208
$ws.Range("M132").Value = -1431.4208
$ws.Range("H136").Value = 4468.1055
$ws.Range("I136").Value = 4159.1816
$ws.Range("K136").Value = 12477.5448
$ws.Range("M136").Value = -9927.5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7585.1665
$ws.Range("I20").Value = 8677.556
$ws.Range("K20").Value = 8677.556
$ws.Range("M20").Value = -8430.556
$ws.Range("H80").Value = 1485.1428
$ws.Range("I80").Value = 1085.6666
$ws.Range("J80").Value = 1594.091
$ws.Range("K80").Value = 1085.6666
$ws.Range("L80").Value = 1594.091
$ws.Range("M80").Value = -87.66660000000002
$ws.Range("N80").Value = -3590.091
$ws.Range("H82").Value = 20063.25
$ws.Range("I82").Value = 13751.333
$ws.Range("J82").Value = 38999
$ws.Range("K82").Value = 13751.333
$ws.Range("L82").Value = 38999
$ws.Range("M82").Value = -13368.333
$ws.Range("N82").Value = -39765
$ws.Range("H83").Value = 1485.1428
$ws.Range("I83").Value = 1085.6666
$ws.Range("J83").Value = 1594.091
$ws.Range("K83").Value = 5428.333000000001
$ws.Range("L83").Value = 7970.455
$ws.Range("M83").Value = -436.3330000000005
$ws.Range("N83").Value = -17954.455
$ws.Range("H85").Value = 20063.25
$ws.Range("I85").Value = 13751.333
$ws.Range("J85").Value = 38999
$ws.Range("K85").Value = 13751.333
$ws.Range("L85").Value = 38999
$ws.Range("M85").Value = -12425.333
$ws.Range("N85").Value = -41651
$ws.Range("H88").Value = 19999.5
$ws.Range("J88").Value = 19999.5
$ws.Range("L88").Value = 19999.5
$ws.Range("N88").Value = -20811.5
$ws.Range("H91").Value = 19999.5
$ws.Range("J91").Value = 19999.5
$ws.Range("L91").Value = 19999.5
$ws.Range("N91").Value = -22807.5
$ws.Range("H105").Value = 1980.8276
$ws.Range("I105").Value = 1585.3334
$ws.Range("K105").Value = 1585.3334
$ws.Range("M105").Value = 161.6666
$ws.Range("H134").Value = 3019.077
$ws.Range("I134").Value = 2549.3635
$ws.Range("K134").Value = 7648.0905
$ws.Range("M134").Value = -5113.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 136.8
$ws.Range("I7").Value = 71
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 71
$ws.Range("L7").Value = 400
$ws.Range("M7").Value = 42
$ws.Range("N7").Value = -626
$ws.Range("H16").Value = 2604.9524
$ws.Range("I16").Value = 2571.8462
$ws.Range("K16").Value = 2571.8462
$ws.Range("M16").Value = -2284.8462
$ws.Range("H22").Value = 799.25
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H52").Value = 66771
$ws.Range("J52").Value = 66771
$ws.Range("L52").Value = 66771
$ws.Range("N52").Value = -67359
$ws.Range("H58").Value = 2982.95
$ws.Range("I58").Value = 3015.647
$ws.Range("J58").Value = 2797.6667
$ws.Range("K58").Value = 3015.647
$ws.Range("L58").Value = 2797.6667
$ws.Range("M58").Value = -2812.647
$ws.Range("N58").Value = -3203.6667
$ws.Range("H113").Value = 2604.9524
$ws.Range("I113").Value = 2571.8462
$ws.Range("K113").Value = 2571.8462
$ws.Range("M113").Value = -401.8462
$ws.Range("H134").Value = 3134.875
$ws.Range("I134").Value = 3431.2222
$ws.Range("K134").Value = 10293.6666
$ws.Range("M134").Value = -7758.6666
$ws.Range("H136").Value = 2982.95
$ws.Range("I136").Value = 3015.647
$ws.Range("J136").Value = 2797.6667
$ws.Range("K136").Value = 9046.940999999999
$ws.Range("L136").Value = 8393.000100000001
$ws.Range("M136").Value = -6496.940999999999
$ws.Range("N136").Value = -13493.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6332.3335
$ws.Range("I39").Value = 179.5
$ws.Range("K39").Value = 538.5
$ws.Range("M39").Value = -244.5
$ws.Range("H74").Value = 13333.333
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -13939
$ws.Range("H77").Value = 13333.333
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -39696
$ws.Range("H140").Value = 1309.8
$ws.Range("I140").Value = 1309.8
$ws.Range("K140").Value = 3929.4
$ws.Range("M140").Value = 1250.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3390.926
$ws.Range("I102").Value = 3630
$ws.Range("J102").Value = 3133.4614
$ws.Range("K102").Value = 3630
$ws.Range("L102").Value = 3133.4614
$ws.Range("M102").Value = -2008
$ws.Range("N102").Value = -6377.4614
$ws.Range("H132").Value = 4221.6665
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 3999.2856
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 11997.8568
$ws.Range("M132").Value = -12470
$ws.Range("N132").Value = -17057.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 5101.8
$ws.Range("J19").Value = 25000
$ws.Range("L19").Value = 25000
$ws.Range("N19").Value = -25340
$ws.Range("H46").Value = 2489.9375
$ws.Range("I46").Value = 2494.75
$ws.Range("J46").Value = 2488.3333
$ws.Range("K46").Value = 2494.75
$ws.Range("L46").Value = 2488.3333
$ws.Range("M46").Value = -2306.75
$ws.Range("N46").Value = -2864.3333
$ws.Range("H55").Value = 305.04544
$ws.Range("I55").Value = 373.7
$ws.Range("K55").Value = 373.7
$ws.Range("M55").Value = -200.7
$ws.Range("H109").Value = 76998.5
$ws.Range("J109").Value = 76998.5
$ws.Range("L109").Value = 76998.5
$ws.Range("N109").Value = -79772.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4279.7646
$ws.Range("I132").Value = 4649.593
$ws.Range("K132").Value = 13948.779
$ws.Range("M132").Value = -11418.779
$ws.Range("H136").Value = 5002999.5
$ws.Range("I136").Value = 5002999.5
$ws.Range("K136").Value = 15008998.5
$ws.Range("M136").Value = -15006448.5
